$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 109; this shifts existing rows 109..142 down to 110..143
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new weekly record
$ws.Cells.Item(109, 1).Value = 4
$ws.Cells.Item(109, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(109, 3).Value = "Los Lagos"
$ws.Cells.Item(109, 4).Value = 44463
$ws.Cells.Item(109, 5).Value = 10
$ws.Cells.Item(109, 6).Value = "Fruta"
$ws.Cells.Item(109, 7).Value = 100102
$ws.Cells.Item(109, 8).Value = "Cítricos"
$ws.Cells.Item(109, 9).Value = 100102006
$ws.Cells.Item(109, 10).Value = "Pomelo"
$ws.Cells.Item(109, 11).Value = "Start Ruby"
$ws.Cells.Item(109, 12).Value = "Primera"
$ws.Cells.Item(109, 13).Value = 160
$ws.Cells.Item(109, 14).Value = 12000
$ws.Cells.Item(109, 15).Value = 12000
$ws.Cells.Item(109, 16).Value = 12000
$ws.Cells.Item(109, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(109, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(109, 19).Value = 857
$ws.Cells.Item(109, 20).Value = 14
